# Apply "Q5/W5" quiz+weekly-assignment columns to the gradebook sheet,
# fix the mis-scanned "?" entries in row 23 (Sturdifen, Jasmine), and
# move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new Q5 / W5 columns -------------------------------------
$ws.Range("M1").Value2 = "Q5"
$ws.Range("N1").Value2 = "W5"

# --- Fix row 23 (Sturdifen, Jasmine): replace placeholder "?" text with
#     the actual numeric Q3/W3 scores -------------------------------------
$ws.Range("E23").Value2 = 4
$ws.Range("F23").Value2 = 48

# --- New Q5 / W5 scores for every student row (2-25) ----------------------
$q5 = @{
    2  = 2;  3  = 4;  4  = 4;  5  = 2;  6  = 0;  7  = 4;  8  = 0;
    9  = 4;  10 = 2;  11 = 2;  12 = 4;  13 = 2;  14 = 4;  15 = 0;
    16 = 2;  17 = 2;  18 = 2;  19 = 4;  20 = 0;  21 = 2;  22 = 2;
    23 = 4;  24 = 0;  25 = 2
}
$w5 = @{
    2  = 50; 3  = 38; 4  = 48; 5  = 45; 6  = 0;  7  = 48; 8  = 0;
    9  = 48; 10 = 50; 11 = 38; 12 = 43; 13 = 39; 14 = 45; 15 = 45;
    16 = 43; 17 = 50; 18 = 50; 19 = 39; 20 = 39; 21 = 38; 22 = 39;
    23 = 38; 24 = 43; 25 = 45
}

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 13).Value2 = $q5[$r]
    $ws.Cells.Item($r, 14).Value2 = $w5[$r]
}

# --- Move the active selection to M8, as in the saved workbook -----------
$ws.Range("M8").Select() | Out-Null
